# Swap the "Material" / "Quantity" / "Inventoryitem" (columns D, E, H) data
# between the "Black" ink row (row 3) and the "Varnish" ink row (row 5), and
# between the "PANTONE Yellow U" ink row (row 8) and the "Magenta" ink row
# (row 9). This reproduces the reordering of those ink rows while leaving
# every other row/column (and all cell styles) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force a value to be written back as text (shared
# string) rather than being auto-coerced into a number (e.g. "1.05"),
# and to let PasteSpecial carry just the value onto the destination so the
# destination's own style/format is left completely untouched.
$scratchAddr = "ZZ999"

function Set-TextValue($destRange, [string]$text) {
    if ($text -eq "") {
        # An empty string can't round-trip through the scratch-cell/
        # PasteSpecial dance below (a blank cell carries nothing to copy),
        # so just clear the destination directly.
        $destRange.Value = ""
        return
    }
    $scratch = $ws.Range($scratchAddr)
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $destRange.PasteSpecial(-4163)  # xlPasteValues
    $scratch.ClearContents()
}

function Swap-Rows($row1, $row2, $cols) {
    # Read every "before" value first so none of them get clobbered while
    # we're still writing the others.
    $before1 = @{}
    $before2 = @{}
    foreach ($col in $cols) {
        $before1[$col] = $ws.Range("$col$row1").Text
        $before2[$col] = $ws.Range("$col$row2").Text
    }
    foreach ($col in $cols) {
        Set-TextValue ($ws.Range("$col$row1")) $before2[$col]
        Set-TextValue ($ws.Range("$col$row2")) $before1[$col]
    }
}

# Columns D (Material), E (Quantity), H (Inventoryitem)
Swap-Rows 3 5 @("D", "E", "H")
Swap-Rows 8 9 @("D", "E", "H")

# Drop the scratch row entirely so the used range / dimension isn't
# permanently widened by our bookkeeping cell.
$ws.Range($scratchAddr).EntireRow.Delete()
